# FPL 2023-24 - GW11 Helicopter View: convert the GW-label "ladder" formulas
# in row 4 to their plain values, unhide column A, and move the frozen-pane
# view / active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Replace the chained "GW"&RIGHT(...)+1 formulas in row 4 with their
#        literal text values (GW2 .. GW38), leaving the non-formula cells
#        (C4=GW1, N4/T4/X4 "EU" markers, etc.) untouched.
$gwCells = @(
    "D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","O4","S4",
    "U4","V4","W4","Y4","AA4","AB4","AC4","AG4","AL4","AM4","AO4",
    "AQ4","AS4","AU4","AW4","AY4","BC4","BD4","BE4","BG4","BI4",
    "BK4","BM4","BO4","BQ4"
)
$gwValues = @(
    "GW2","GW3","GW4","GW5","GW6","GW7","GW8","GW9","GW10","GW11","GW12","GW13",
    "GW14","GW15","GW16","GW17","GW18","GW19","GW20","GW21","GW22","GW23","GW24",
    "GW25","GW26","GW27","GW28","GW29","GW30","GW31","GW32","GW33","GW34",
    "GW35","GW36","GW37","GW38"
)

for ($i = 0; $i -lt $gwCells.Length; $i++) {
    $ws.Range($gwCells[$i]).Value = $gwValues[$i]
}

# --- 2) Unhide column A (was hidden="1").
$ws.Columns.Item(1).Hidden = $false

# --- 3) Update the frozen pane top-left cell and the active selection.
$ws.Range("C1").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("P32").Select()
